# Rode draad nr. 10: adjust bepaal_N_vak.
# Old:  N_vak = 1 + a * lengte / delta_L
# New:  N_vak = MAX(1, a * lengte / delta_L)
# pf_vak is recomputed from the (possibly updated) N_vak: pf_vak = N_vak * pf_dsn

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colLengte = 4   # D: lengte
$colA      = 8   # H: a
$colPfDsn  = 9   # I: pf_dsn
$colDeltaL = 12  # L: delta_L
$colNVak   = 14  # N: N_vak
$colPfVak  = 15  # O: pf_vak

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($row = 2; $row -le $lastRow; $row++) {
    $a      = $ws.Cells.Item($row, $colA).Value2
    $lengte = $ws.Cells.Item($row, $colLengte).Value2
    $deltaL = $ws.Cells.Item($row, $colDeltaL).Value2
    $pfDsn  = $ws.Cells.Item($row, $colPfDsn).Value2

    $raw = $a * $lengte / $deltaL

    if ($raw -lt 1) {
        $nVak = 1
    } else {
        $nVak = $raw
    }

    $pfVak = $nVak * $pfDsn

    $ws.Cells.Item($row, $colNVak).Value2 = $nVak
    $ws.Cells.Item($row, $colPfVak).Value2 = $pfVak
}
